$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "moving" data block (A,B,E,F,G,H,K-N presence,Q,R,AC) for rows 2-13
# based on a full permutation of per-observation records across rows.

# Row 2  (Id 111815508)
$ws.Range("A2").Value = 111815508
$ws.Range("B2").Value = 56398
$ws.Range("E2").Value = 100109
$ws.Range("F2").Value = 'Tretåig hackspett'
$ws.Range("G2").Value = 'Picoides tridactylus'
$ws.Range("H2").Value = '(Linnaeus, 1758)'
$ws.Range("Q2").Value = 458162.4570845839
$ws.Range("R2").Value = 7054329.489790585
$ws.Range("I2").Copy($ws.Range("K2"))
$ws.Range("I2").Copy($ws.Range("L2"))
$ws.Range("I2").Copy($ws.Range("M2"))
$ws.Range("I2").Copy($ws.Range("N2"))
$ws.Range("AC2").Value = 'ringhack'

# Row 3  (Id 111815518)
$ws.Range("A3").Value = 111815518
$ws.Range("B3").Value = 77515
$ws.Range("E3").Value = 6425
$ws.Range("F3").Value = 'Garnlav'
$ws.Range("G3").Value = 'Alectoria sarmentosa'
$ws.Range("H3").Value = '(Ach.) Ach.'
$ws.Range("Q3").Value = 458250.901553072
$ws.Range("R3").Value = 7054618.376188213
$ws.Range("K3").ClearContents()
$ws.Range("L3").ClearContents()
$ws.Range("M3").ClearContents()
$ws.Range("N3").ClearContents()
$ws.Range("AC3").ClearContents()

# Row 4  (Id 111815516)
$ws.Range("A4").Value = 111815516
$ws.Range("B4").Value = 89423
$ws.Range("E4").Value = 5432
$ws.Range("F4").Value = 'Granticka'
$ws.Range("G4").Value = 'Porodaedalea chrysoloma'
$ws.Range("H4").Value = '(Fr.) Fiasson & Niemelä'
$ws.Range("Q4").Value = 458289.5512131723
$ws.Range("R4").Value = 7054475.069158822
$ws.Range("K4").ClearContents()
$ws.Range("L4").ClearContents()
$ws.Range("M4").ClearContents()
$ws.Range("N4").ClearContents()
$ws.Range("AC4").ClearContents()

# Row 5  (Id 111815512)
$ws.Range("A5").Value = 111815512
$ws.Range("B5").Value = 56398
$ws.Range("E5").Value = 100109
$ws.Range("F5").Value = 'Tretåig hackspett'
$ws.Range("G5").Value = 'Picoides tridactylus'
$ws.Range("H5").Value = '(Linnaeus, 1758)'
$ws.Range("Q5").Value = 458154.6107204149
$ws.Range("R5").Value = 7054646.336103803
$ws.Range("I2").Copy($ws.Range("K5"))
$ws.Range("I2").Copy($ws.Range("L5"))
$ws.Range("I2").Copy($ws.Range("M5"))
$ws.Range("I2").Copy($ws.Range("N5"))
$ws.Range("AC5").Value = 'ringhack'

# Row 6  (Id 111815517)
$ws.Range("A6").Value = 111815517
$ws.Range("B6").Value = 77515
$ws.Range("E6").Value = 6425
$ws.Range("F6").Value = 'Garnlav'
$ws.Range("G6").Value = 'Alectoria sarmentosa'
$ws.Range("H6").Value = '(Ach.) Ach.'
$ws.Range("Q6").Value = 458250.8216980004
$ws.Range("R6").Value = 7054375.482693202
$ws.Range("K6").ClearContents()
$ws.Range("L6").ClearContents()
$ws.Range("M6").ClearContents()
$ws.Range("N6").ClearContents()
$ws.Range("AC6").ClearContents()

# Row 7  (Id 111815519)
$ws.Range("A7").Value = 111815519
$ws.Range("B7").Value = 77515
$ws.Range("E7").Value = 6425
$ws.Range("F7").Value = 'Garnlav'
$ws.Range("G7").Value = 'Alectoria sarmentosa'
$ws.Range("H7").Value = '(Ach.) Ach.'
$ws.Range("Q7").Value = 458215.7474518137
$ws.Range("R7").Value = 7054621.063481365
$ws.Range("K7").ClearContents()
$ws.Range("L7").ClearContents()
$ws.Range("M7").ClearContents()
$ws.Range("N7").ClearContents()
$ws.Range("AC7").ClearContents()

# Row 8  (Id 111815513)
$ws.Range("A8").Value = 111815513
$ws.Range("B8").Value = 56398
$ws.Range("E8").Value = 100109
$ws.Range("F8").Value = 'Tretåig hackspett'
$ws.Range("G8").Value = 'Picoides tridactylus'
$ws.Range("H8").Value = '(Linnaeus, 1758)'
$ws.Range("Q8").Value = 458173.7327805056
$ws.Range("R8").Value = 7054711.474791372
$ws.Range("I2").Copy($ws.Range("K8"))
$ws.Range("I2").Copy($ws.Range("L8"))
$ws.Range("I2").Copy($ws.Range("M8"))
$ws.Range("I2").Copy($ws.Range("N8"))
$ws.Range("AC8").Value = 'ringhack gamla'

# Row 9  (Id 111815515)
$ws.Range("A9").Value = 111815515
$ws.Range("B9").Value = 89423
$ws.Range("E9").Value = 5432
$ws.Range("F9").Value = 'Granticka'
$ws.Range("G9").Value = 'Porodaedalea chrysoloma'
$ws.Range("H9").Value = '(Fr.) Fiasson & Niemelä'
$ws.Range("Q9").Value = 458161.9437607233
$ws.Range("R9").Value = 7054459.400503729
$ws.Range("K9").ClearContents()
$ws.Range("L9").ClearContents()
$ws.Range("M9").ClearContents()
$ws.Range("N9").ClearContents()
$ws.Range("AC9").ClearContents()

# Row 10  (Id 111815510)
$ws.Range("A10").Value = 111815510
$ws.Range("B10").Value = 56398
$ws.Range("E10").Value = 100109
$ws.Range("F10").Value = 'Tretåig hackspett'
$ws.Range("G10").Value = 'Picoides tridactylus'
$ws.Range("H10").Value = '(Linnaeus, 1758)'
$ws.Range("Q10").Value = 458203.7272220219
$ws.Range("R10").Value = 7054385.000644128
$ws.Range("I2").Copy($ws.Range("K10"))
$ws.Range("I2").Copy($ws.Range("L10"))
$ws.Range("I2").Copy($ws.Range("M10"))
$ws.Range("I2").Copy($ws.Range("N10"))
$ws.Range("AC10").Value = 'ringhack'

# Row 11  (Id 111815514)
$ws.Range("A11").Value = 111815514
$ws.Range("B11").Value = 89423
$ws.Range("E11").Value = 5432
$ws.Range("F11").Value = 'Granticka'
$ws.Range("G11").Value = 'Porodaedalea chrysoloma'
$ws.Range("H11").Value = '(Fr.) Fiasson & Niemelä'
$ws.Range("Q11").Value = 458153.7808649908
$ws.Range("R11").Value = 7054482.19637617
$ws.Range("K11").ClearContents()
$ws.Range("L11").ClearContents()
$ws.Range("M11").ClearContents()
$ws.Range("N11").ClearContents()
$ws.Range("AC11").ClearContents()

# Row 12  (Id 111815509)
$ws.Range("A12").Value = 111815509
$ws.Range("B12").Value = 56398
$ws.Range("E12").Value = 100109
$ws.Range("F12").Value = 'Tretåig hackspett'
$ws.Range("G12").Value = 'Picoides tridactylus'
$ws.Range("H12").Value = '(Linnaeus, 1758)'
$ws.Range("Q12").Value = 458176.2590895323
$ws.Range("R12").Value = 7054362.673967168
$ws.Range("I2").Copy($ws.Range("K12"))
$ws.Range("I2").Copy($ws.Range("L12"))
$ws.Range("I2").Copy($ws.Range("M12"))
$ws.Range("I2").Copy($ws.Range("N12"))
$ws.Range("AC12").Value = 'ringhack gamla'

# Row 13  (Id 111815507)
$ws.Range("A13").Value = 111815507
$ws.Range("B13").Value = 56398
$ws.Range("E13").Value = 100109
$ws.Range("F13").Value = 'Tretåig hackspett'
$ws.Range("G13").Value = 'Picoides tridactylus'
$ws.Range("H13").Value = '(Linnaeus, 1758)'
$ws.Range("Q13").Value = 458151.5539710881
$ws.Range("R13").Value = 7054482.225765129
$ws.Range("I2").Copy($ws.Range("K13"))
$ws.Range("I2").Copy($ws.Range("L13"))
$ws.Range("I2").Copy($ws.Range("M13"))
$ws.Range("I2").Copy($ws.Range("N13"))
$ws.Range("AC13").Value = 'ringhack gamla'
